$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 11-13 (col G) as Done
$ws.Range("G11").Value = "Done"
$ws.Range("G12").Value = "Done"
$ws.Range("G13").Value = "Done"

# Append new schedule rows 35-46
$ws.Range("D35").Value = 44907
$ws.Range("E35").Value = "https://zoom.us/rec/share/U-GB6_oep1BC7i1Pi_YZW96fh6u3ZtCm2vqYghfLkHT7T4X1tg6mkJKH6C0sJ4yx.BooSWKSJUtkO_FnK"
$ws.Range("D36").Value = 44909
$ws.Range("E36").Value = "https://zoom.us/rec/share/Psrs6H-mVO7Vf_tbnlVJbsffZs5vuL1phykJSYkR35uZvLvqOJdOA2Ydf3dcXY4p.OtLAMDM3tEY5jE13"
$ws.Range("D37").Value = 44911
$ws.Range("E37").Value = "https://zoom.us/rec/share/9k79XAVdCgQEtyN-gHK0aElB-eS9wxbHniJpqUtaiDiZKd2p2HE1LYETsjlJuHm9.2zELiyX7ZuG1BiMR"
$ws.Range("D38").Value = 44914
$ws.Range("E38").Value = "https://zoom.us/rec/share/40jbucMo5iPRva_a2X6fNEnz-z5D2yWjsyXMdfrjgbXW2Mft8Va1dV9aszQAQ-BV.Z2g9Ywbj-zSgC52D"
$ws.Range("D39").Value = 44916
$ws.Range("E39").Value = "https://zoom.us/rec/share/S6_IqmA2fQA8vb-AF8c3X3nwAgEZUL_yCWWTsyI_blAqMhYU9tyxnIazMBUv7u6N.mMh9RREU9nxf4Jk4"
$ws.Range("D40").Value = 44918
$ws.Range("E40").Value = "Nghỉ??"
$ws.Range("D41").Value = 44921
$ws.Range("E41").Value = "https://zoom.us/rec/share/-Ryk4BEGQNEkjKkeGA3q_8HUT-XE8F76UdGU1y2N0luxB4XKmBQa60oi0Rql4Nr1.bIwqmW_DQLI6fW9A"
$ws.Range("D42").Value = 44923
$ws.Range("E42").Value = "https://zoom.us/rec/share/x8QSnXuJoJn5OEtSVVtK6hdWcyMeuckQSMNTAzCwkXApN7yBkzk0jGZNk9lW2USJ.p2n9T7RCiC6eQRSQ"
$ws.Range("D43").Value = 44925
$ws.Range("E43").Value = "Nghỉ??"
$ws.Range("D44").Value = 44928
$ws.Range("E44").Value = "Nghỉ??"
$ws.Range("D45").Value = 44930
$ws.Range("E45").Value = "https://zoom.us/rec/share/CQ35g_0cFAOENXmBJcM_XwMdhWc_XDbxxTMuIoe-UlhD4_WjhDV0jgN4N66Kz39W.Rz0rmKwKEp5BonWA"
$ws.Range("D46").Value = 44932

# Freeze header row/columns and set view selection to match the author's final state
$ws.Range("D4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E42").Select()
